$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68..121 down to 69..122.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new price record.
$ws.Cells.Item(68, 1).Value2 = 7
$ws.Cells.Item(68, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value2 = "Ñuble"
$ws.Cells.Item(68, 4).Value2 = 44574
$ws.Cells.Item(68, 5).Value2 = 16
$ws.Cells.Item(68, 6).Value2 = 100112045
$ws.Cells.Item(68, 7).Value2 = "Zapallo"
$ws.Cells.Item(68, 8).Value2 = "Camote"
$ws.Cells.Item(68, 9).Value2 = "1a nueva(o)"
$ws.Cells.Item(68, 10).Value2 = 200
$ws.Cells.Item(68, 11).Value2 = 300
$ws.Cells.Item(68, 12).Value2 = 350
$ws.Cells.Item(68, 13).Value2 = 325
$ws.Cells.Item(68, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(68, 15).Value2 = "Región del Maule"
$ws.Cells.Item(68, 16).Value2 = 325
$ws.Cells.Item(68, 17).Value2 = 1
$ws.Cells.Item(68, 18).Value2 = "Hortaliza"
